# Update the cryptocurrency price/volume table with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.915.86"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.816.17"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'309.15"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4637"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").Value = "'0.3661"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").Value = "'0.07359"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "'0.8717"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "'20.30"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.817.78"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "'5.373"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'0.07099"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "'6.505"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "'91.56"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'0.000008722"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'14.66"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "26.953.86"
$ws.Range("D22").Value = "'5.300"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'10.63"
$ws.Range("D24").Value = "2.070.00"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'150.87"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'18.32"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'2.135"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'5.258"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'115.43"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'0.08913"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'0.7567"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "'4.481"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'2.908"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'1.002"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "'1.085"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.995"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01951"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'7.236"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "'0.5303"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "'2.312"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "'0.1654"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'8.444"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'0.4869"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "'103.51"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "'1.662"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'0.06292"
$ws.Range("E51").Value = "  -0.01%  "
